$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.249265789985657
$ws.Range("B1").Value = 1.404392957687378
$ws.Range("C1").Value = 1.711446046829224
$ws.Range("D1").Value = 3.251639604568481
$ws.Range("E1").Value = 15
